# Updated cryptos list on Sat Nov  4 17:42:38 UTC 2023 with GitHub Actions
# Refreshes coin Price (D) / Volume(1h) (E) columns; rows 14 & 15 swap coins.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure numeric-looking Price strings stay literal text (not numbers)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.166.79"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.857.55"
$ws.Range("E3").Value = "  +1.96%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "239.11"
$ws.Range("E5").Value = "  +3.63%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.37%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.11%  "

# Row 8 - Solana
$ws.Range("D8").Value = "41.91"
$ws.Range("E8").Value = "  +5.77%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.329"
$ws.Range("E9").Value = "  +3.08%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.66%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.07%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.126.21"
$ws.Range("E12").Value = "  +1.90%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +1.83%  "

# Row 14 - Polygon->WrappedEther
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.848.96"
$ws.Range("E14").Value = "  +0.67%  "

# Row 15 - WrappedEther->Polygon
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.677"
$ws.Range("E15").Value = "  +1.79%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +2.06%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "35.110.04"
$ws.Range("E17").Value = "  +1.32%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "69.73"
$ws.Range("E18").Value = "  +0.39%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +1.45%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "240.61"
$ws.Range("E20").Value = "  +0.66%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  +1.71%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "4.72"
$ws.Range("E22").Value = "  +1.68%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.15%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.04%  "

# Row 25 - Monero
$ws.Range("D25").Value = "168.80"
$ws.Range("E25").Value = "  -2.53%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "1.91"
$ws.Range("E26").Value = "  +27.15%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +3.68%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "17.66"
$ws.Range("E28").Value = "  +2.30%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  +0.33%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.19%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.0557"
$ws.Range("E31").Value = "  +1.87%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.00"
$ws.Range("E32").Value = "  +2.46%  "

# Row 33 - WEMIXToken
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").Value = "  +27.85%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +2.71%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "0.825"
$ws.Range("E35").Value = "  +18.82%  "

# Row 36 - LidoDAOToken
$ws.Range("D36").Value = "2.01"
$ws.Range("E36").Value = "  +10.95%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  +6.81%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +7.62%  "

# Row 39 - Aave
$ws.Range("D39").Value = "90.52"
$ws.Range("E39").Value = "  -0.88%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +4.42%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.342.06"
$ws.Range("E41").Value = "  +0.28%  "

# Row 42 - InjectiveProtocol
$ws.Range("D42").Value = "14.95"
$ws.Range("E42").Value = "  +3.61%  "

# Row 43 - RenderToken
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  +3.99%  "

# Row 44 - HuobiToken
$ws.Range("E44").Value = "  -0.76%  "

# Row 45 - Gas
$ws.Range("D45").Value = "12.37"
$ws.Range("E45").Value = "  +48.26%  "

# Row 46 - Kaspa
$ws.Range("D46").Value = "0.0557"
$ws.Range("E46").Value = "  +6.74%  "

# Row 47 - MXToken
$ws.Range("E47").Value = "  -0.09%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "6.62"
$ws.Range("E48").Value = "  +5.88%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.038.98"
$ws.Range("E49").Value = "  +1.62%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.0680"
$ws.Range("E50").Value = "  +1.88%  "

# Row 51 - PaxDollar
$ws.Range("E51").Value = "  +0.11%  "

